$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "DifferentTypes"

# Header row
$ws2.Range("A1").Value = "StringValue"
$ws2.Range("B1").Value = "NumericValue"
$ws2.Range("C1").Value = "BooleanValue"
$ws2.Range("D1").Value = "FormulaValue"
$ws2.Range("E1").Value = "FormulaErrorValue"

# Data row exercising different cell types
$ws2.Range("A2").Value = "City"
$ws2.Range("B2").Value = 17
$ws2.Range("C2").Value = $false
$ws2.Range("D2").Formula = "=B2*B2"
$ws2.Range("E2").Formula = "=A2*B2"

# Column widths matching the authored layout
$ws2.Range("A1").EntireColumn.ColumnWidth = 10.5
$ws2.Range("B1").EntireColumn.ColumnWidth = 11.666667
$ws2.Range("C1:D1").EntireColumn.ColumnWidth = 13
$ws2.Range("E1").EntireColumn.ColumnWidth = 17.666667

# Selection ends on E2; this also marks the sheet as the active tab
$ws2.Range("E2").Select()
